# Component 5 (C5) PowerPoint edit:
#  1. Slide 6's table is re-styled from the custom "Table_0" style to the
#     built-in "Medium Style 2 - Accent 1" table style.
#  2. The deck's theme (theme1.xml, used by the slide master / all slides)
#     is switched from the "Integral" palette to the stock "Office" palette.
#     (theme2.xml, which only feeds the Notes Master, and the <a:theme>/
#     <a:clrScheme> "name" attributes themselves are not reachable through
#     the PowerPoint object model, so only the 12 scheme colors - the part
#     of the theme that actually drives the deck's look - are updated here.)

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 --------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{8F5A0379-8E82-4B5B-B693-2A6522DC906D}")

# --- 2. Theme colors: "Integral" -> "Office" -------------------------------
# Colors() is 1-based and ordered dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink; RGB is a COLORREF (0xBBGGRR), so each target hex color's bytes
# are reversed.
$themeSlide = $p.Slides.Item(1)
$colorScheme = $themeSlide.ThemeColorScheme

$colorScheme.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$colorScheme.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$colorScheme.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$colorScheme.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$colorScheme.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$colorScheme.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$colorScheme.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$colorScheme.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$colorScheme.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$colorScheme.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$colorScheme.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$colorScheme.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
